$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Data edits (re-randomised condition table, rows 17-21 and 56-61)
# -----------------------------------------------------------------

# Row 17: values unchanged, but D17/E17 pick up the right-aligned
# "Aptos Narrow" style already used by column C (style index 2).
$ws.Range("C17").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("C17").Copy()
$ws.Range("E17").PasteSpecial(-4122)

# Row 18
$ws.Range("B18").Value = "X"
$ws.Range("C17").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("C17").Copy()
$ws.Range("E18").PasteSpecial(-4122)

# Row 19
$ws.Range("B19").Value = "X"
$ws.Range("C19").Value = 1
$ws.Range("C17").Copy()
$ws.Range("D19").PasteSpecial(-4122)
$ws.Range("C17").Copy()
$ws.Range("E19").PasteSpecial(-4122)

# Row 20
$ws.Range("B20").Value = "P"

# Row 21
$ws.Range("B21").Value = "O"
$ws.Range("C21").Value = 0

# Row 56
$ws.Range("B56").Value = "Q"

# Row 57
$ws.Range("B57").Value = "A"
$ws.Range("C57").Value = 0

# Row 60
$ws.Range("B60").Value = "U"
$ws.Range("C60").Copy()
$ws.Range("D60").PasteSpecial(-4122)
$ws.Range("C60").Copy()
$ws.Range("E60").PasteSpecial(-4122)

# Row 61
$ws.Range("B61").Value = "U"
$ws.Range("C61").Value = 1
$ws.Range("C60").Copy()
$ws.Range("D61").PasteSpecial(-4122)
$ws.Range("C60").Copy()
$ws.Range("E61").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# -----------------------------------------------------------------
# View state: selected cell moves from A2:C61 (active A2) to F4
# -----------------------------------------------------------------
$ws.Range("F4").Select()

# Best-effort: the workbook window had scrolled/moved on screen
# (xWindow 2000 -> 10960) before the file was re-saved.
$excel.ActiveWindow.Left = 10960
